$d = $word.ActiveDocument

# Locate the exact sentence that needs to be split into three runs.
$full = $d.Content
$found = $full.Find.Execute(
    "Add current sources and different types of loads and converters, modify the file schema",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Anchor on the boundary right after "loads " (before "and converters").
    $boundary1 = $full.Duplicate
    $boundary1.Find.Execute("loads ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $boundary1.Collapse(0)

    # Insert the new parenthetical run text at that boundary.
    $boundary1.InsertAfter("(constant power load, E-lock and OBC board) ")

    # Toggling formatting on the freshly inserted text forces it to stay a
    # distinct run instead of being re-coalesced into its neighbours, even
    # though the net formatting ends up identical to the surrounding runs.
    $boundary1.Bold = 1
    $boundary1.Bold = 0
}
